# NIT-8600790247.xlsx — "Elimina EC anteriores y se agregan nuevos,
# se modifica base de datos"
#
# The mora-detail table (rows 16-17, columns E "Periodo Mora" / F "Valor
# Mora") gets its two records swapped:
#   Row 16 (before): Periodo 2502, Valor  189000
#   Row 17 (before): Periodo 2503, Valor    6300
#   Row 16 (after) : Periodo 2503, Valor    6300
#   Row 17 (after) : Periodo 2502, Valor  189000
# All other cells/styles/rows are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E16").Value = "2503"
$ws.Range("F16").Value = 6300

$ws.Range("E17").Value = "2502"
$ws.Range("F17").Value = 189000
